$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 (main schedule sheet) ---
$ws1 = $wb.Worksheets.Item("LP1912")

# A new scrape row needs to be inserted between the existing row 110 and
# what is currently row 111 (shifts old rows 111-118 down to 112-119).
$ws1.Rows.Item(111).Insert()

# Fill in the newly inserted row 111.
$ws1.Range("A111").Value = "08:48:08"
$ws1.Range("B111").Value = "09:35"
$ws1.Range("C111").Value = "23_HERNANDEZ"
$ws1.Range("D111").Value = 47
$ws1.Range("E111").Value = "LP1912"

# Append two brand new rows (120 and 121) at the end of the table.
$ws1.Range("A120").Value = "08:48:08"
$ws1.Range("B120").Value = "10:42"
$ws1.Range("C120").Value = "17_ROMERO"
$ws1.Range("D120").Value = 114
$ws1.Range("E120").Value = "LP1912"

$ws1.Range("A121").Value = "08:48:08"
$ws1.Range("B121").Value = "10:44"
$ws1.Range("C121").Value = "14_ABASTO"
$ws1.Range("D121").Value = 116
$ws1.Range("E121").Value = "LP1912"

# Refresh the scrape metadata header on the main sheet.
$ws1.Range("A2").Value = "Última actualización: 08:48:08"
$ws1.Range("A3").Value = "Total filas: 116"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:48:08"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:48:08"
